# Update forest data - 2026-01-28 12:23
#
# The "New" sheet's 5 current listings age out and move down into the
# "Previously added" sheet (appended as new rows), and 2 fresh listings
# take their place on "New".

$wb  = $excel.ActiveWorkbook
$old = $wb.Worksheets.Item("Previously added")
$new = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------
# 1) Capture the 5 existing hyperlink targets on "New" (rows 2..6) before
#    anything is touched, so we can re-attach them on "Previously added".
# ---------------------------------------------------------------------
$links = @()
for ($r = 2; $r -le 6; $r++) {
    $links += $new.Range("A" + $r).Hyperlinks.Item(1).Address
}

# ---------------------------------------------------------------------
# 2) Move the 5 rows (values + number formats + styles) down onto
#    "Previously added", right after the existing data (row 434).
# ---------------------------------------------------------------------
$destFirstRow = $old.UsedRange.Rows.Count + 1   # 435
$new.Range("A2:F6").Copy($old.Range("A" + $destFirstRow))

for ($i = 0; $i -lt 5; $i++) {
    $destRow = $destFirstRow + $i
    $old.Hyperlinks.Add($old.Range("A" + $destRow), $links[$i])
}

# ---------------------------------------------------------------------
# 3) Reset "New": drop its hyperlinks, drop the 3 rows that are no
#    longer needed (it will only hold 2 fresh listings), then fill in
#    the 2 new rows.
# ---------------------------------------------------------------------
$new.Range("A2").Hyperlinks.Delete()
$new.Rows("4:6").Delete()

$new.Range("A2").Value = "https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/jersikas-pag/fpncc.html"
$new.Range("B2").Value = "8 500 €"
$new.Range("C2").Value = "Preiļi un raj."
$new.Range("D2").Value = "2.50 ha."
$new.Range("E2").Value = "7652 004 0129"
$new.Range("F2").Value = 46050.56597222222

$new.Range("A3").Value = "https://www.ss.com/msg/lv/real-estate/wood/saldus-and-reg/zvardes-pag/kjfkn.html"
$new.Range("B3").Value = "89 000 €"
$new.Range("C3").Value = "Saldus un raj."
$new.Range("D3").Value = "16 ha."
$new.Range("E3").Formula = '=TEXT(84980020059,"0")'
$new.Range("E3").Copy()
$new.Range("E3").PasteSpecial(-4163)   # xlPasteValues - keep the text, drop the formula
$new.Range("F3").Value = 46050.527083333334

$new.Hyperlinks.Add($new.Range("A2"), "https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/jersikas-pag/fpncc.html")
$new.Hyperlinks.Add($new.Range("A3"), "https://www.ss.com/msg/lv/real-estate/wood/saldus-and-reg/zvardes-pag/kjfkn.html")
